$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (cleared)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# (write via a scratch cell forced to Text format so the ISO-looking string
#  isn't auto-converted into a date serial, then paste-special the *value*
#  back onto B8 so the destination keeps its original General style)
$scratch = $ws.Range("Z1000")
$scratch.NumberFormat = "@"
$scratch.Value = "2025-11-18"
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)
$scratch.Clear()
